# Update workbook data: rotate the date/volume/price fields among rows 2-4
# so that row2 <- old row3, row3 <- old row4, row4 <- old row2 (for columns
# D, M, N, O, P, Q, S, T), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was 44875 / 50 / 16000 / 16000 / 16000 / "$/bandeja 10 kilos" / 1600 / 10)
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 67

# Row 3 (was 44874 / 67 / 16000 / 16000 / 16000 / "$/bandeja 10 kilos" / 1600 / 10)
$ws.Range("D3").Value = 44855
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("S3").Value = 3000
$ws.Range("T3").Value = 5

# Row 4 (was 44855 / 25 / 15000 / 15000 / 15000 / "$/bandeja 5 kilos" / 3000 / 5)
$ws.Range("D4").Value = 44875
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1600
$ws.Range("T4").Value = 10
